# HRMViewer: show cluster/coordinator images
# Group each "Ellipse N" + its associated "Textfeld N" label shape into a
# single p:grpSp, one group per cluster-level indicator (levels 0-4).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Level 0: "Ellipse 3" (shape 1) + "Textfeld 4" (shape 2) -> "Gruppieren 1"
$grp1 = $s.Shapes.Range(@(1, 2)).Group()
$grp1.Name = "Gruppieren 1"

# Level 1: "Ellipse 5" + "Textfeld 6" -> "Gruppieren 2"
$grp2 = $s.Shapes.Range(@(2, 3)).Group()
$grp2.Name = "Gruppieren 2"

# Level 2: "Ellipse 7" + "Textfeld 8" -> "Gruppieren 13"
$grp3 = $s.Shapes.Range(@(3, 4)).Group()
$grp3.Name = "Gruppieren 13"

# Level 3: "Ellipse 9" + "Textfeld 10" -> "Gruppieren 14"
$grp4 = $s.Shapes.Range(@(4, 5)).Group()
$grp4.Name = "Gruppieren 14"

# Level 4: "Ellipse 11" + "Textfeld 12" -> "Gruppieren 15"
$grp5 = $s.Shapes.Range(@(5, 6)).Group()
$grp5.Name = "Gruppieren 15"
